$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold plain text in this sheet (values
# like "1.676.98" or "  +0.21%  " are not real numbers). A handful of the
# refreshed prices (e.g. "215.18") look like genuine decimal numbers, so a
# leading apostrophe is used for those to force Excel to keep them as text,
# exactly like a user typing '215.18 into the cell.
$cellUpdates = [ordered]@{
    "D2" = "26.973.42"
    "E2" = "  +0.06%  "
    "D3" = "1.679.69"
    "E3" = "  +0.42%  "
    "E4" = "  -0.07%  "
    "D5" = "'215.18"
    "E5" = "  -0.39%  "
    "E6" = "  -2.77%  "
    "E7" = "  -0.11%  "
    "D8" = "'0.251"
    "E8" = "  -1.09%  "
    "D9" = "'0.0620"
    "E9" = "  -0.24%  "
    "D10" = "'20.81"
    "E10" = "  +3.01%  "
    "E11" = "  -0.31%  "
    "D12" = "1.913.32"
    "E12" = "  +0.24%  "
    "D13" = "1.679.85"
    "E13" = "  +0.50%  "
    "E14" = "  +0.36%  "
    "D15" = "'0.532"
    "E15" = "  +1.98%  "
    "D16" = "'65.83"
    "E16" = "  +0.19%  "
    "D17" = "26.981.24"
    "E17" = "  +0.01%  "
    "D19" = "'236.50"
    "E19" = "  +1.69%  "
    "D20" = "0.0₃0735"
    "E20" = "  -0.05%  "
    "E21" = "  -0.07%  "
    "D22" = "'4.44"
    "E22" = "  -0.42%  "
    "D23" = "'9.21"
    "E24" = "  -4.12%  "
    "D25" = "'146.45"
    "E25" = "  +0.57%  "
    "D26" = "'7.23"
    "E26" = "  +1.24%  "
    "D27" = "'16.06"
    "E27" = "  +1.31%  "
    "E29" = "  -0.01%  "
    "E30" = "  +0.59%  "
    "E31" = "  -0.59%  "
    "D32" = "'3.34"
    "E32" = "  +0.61%  "
    "D33" = "1.490.42"
    "D34" = "'3.18"
    "E34" = "  +1.06%  "
    "D35" = "'1.69"
    "E35" = "  +5.08%  "
    "E36" = "  +0.04%  "
    "D37" = "'0.585"
    "E37" = "  +3.38%  "
    "D38" = "'0.916"
    "E38" = "  +1.47%  "
    "D39" = "'0.0175"
    "E39" = "  +3.63%  "
    "E40" = "  +6.94%  "
    "E41" = "  -4.98%  "
    "E42" = "  -0.07%  "
    "D43" = "'67.45"
    "E43" = "  +2.44%  "
    "D44" = "'2.27"
    "E44" = "  -1.62%  "
    "D45" = "1.818.93"
    "E45" = "  +0.37%  "
    "D46" = "'0.781"
    "E46" = "  -0.01%  "
    "D47" = "'90.62"
    "E47" = "  -0.11%  "
}

foreach ($ref in $cellUpdates.Keys) {
    $ws.Range($ref).Value = $cellUpdates[$ref]
}

# A new coin (BabyDogeCoin) entered the top-50 list, pushing RenderToken,
# Algorand and EnergySwap down one row each; Cronos falls out of the list.
$rowUpdates = [ordered]@{
    48 = @("BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.0₆0105", "  +2.39%  ")
    49 = @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "'1.52", "  -0.58%  ")
    50 = @("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "'0.103", "  +2.66%  ")
    51 = @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "'7.77", "  +2.11%  ")
}

foreach ($r in $rowUpdates.Keys) {
    $vals = $rowUpdates[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
}
